$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new leaderboard rows (57-69). Row 56 is used as a formatting template
# for "plain" rows (style indices 2/3/3/4/3 on A/B/C/D/E), and D38 (which already
# carries the true "Hyperlink" cell style + border) is used as the formatting
# template for the two rows whose D column becomes a real hyperlink (62, 66).

# Row 57
$ws.Range("A56:F56").Copy()
$ws.Range("A57:F57").PasteSpecial(-4122)
$ws.Range("A57").Value = 45384.081643518519
$ws.Range("B57").Value = "Sachin Pant"
$ws.Range("C57").Value = "B23229"
$ws.Range("D57").Value = "https://www.beecrowd.com.br/judge/en/profile/948971"
$ws.Range("E57").Value = "DSE"
$ws.Range("F57").Value = 0
$ws.Rows.Item(57).RowHeight = 27.6

# Row 58
$ws.Range("A56:F56").Copy()
$ws.Range("A58:F58").PasteSpecial(-4122)
$ws.Range("A58").Value = 45384.101875
$ws.Range("B58").Value = "ADITYA TAYAL"
$ws.Range("C58").Value = "B23243"
$ws.Range("D58").Value = "https://www.beecrowd.com.br/judge/en/profile/948984"
$ws.Range("E58").Value = "EE"
$ws.Range("F58").Value = 0
$ws.Rows.Item(58).RowHeight = 27.6

# Row 59
$ws.Range("A56:F56").Copy()
$ws.Range("A59:F59").PasteSpecial(-4122)
$ws.Range("A59").Value = 45384.167141203703
$ws.Range("B59").Value = "Nitin chaurasiya"
$ws.Range("C59").Value = "B23220"
$ws.Range("D59").Value = "https://www.beecrowd.com.br/judge/en/profile/948196"
$ws.Range("E59").Value = "DSE"
$ws.Range("F59").Value = 0
$ws.Rows.Item(59).RowHeight = 27.6

# Row 60
$ws.Range("A56:F56").Copy()
$ws.Range("A60:F60").PasteSpecial(-4122)
$ws.Range("A60").Value = 45384.188750000001
$ws.Range("B60").Value = "Nitin Chaurasiya"
$ws.Range("C60").Value = "B23220"
$ws.Range("D60").Value = "https://www.beecrowd.com.br/judge/en/profile/948196"
$ws.Range("E60").Value = "DSE"
$ws.Range("F60").Value = 0
$ws.Rows.Item(60).RowHeight = 27.6

# Row 61
$ws.Range("A56:F56").Copy()
$ws.Range("A61:F61").PasteSpecial(-4122)
$ws.Range("A61").Value = 45384.424675925926
$ws.Range("B61").Value = "Harshit Kumar Singh"
$ws.Range("C61").Value = "B23133"
$ws.Range("D61").Value = "https://www.beecrowd.com.br/judge/en/profile/949109"
$ws.Range("E61").Value = "CSE"
$ws.Range("F61").Value = 0
$ws.Rows.Item(61).RowHeight = 27.6

# Row 62
$ws.Range("A56:F56").Copy()
$ws.Range("A62:F62").PasteSpecial(-4122)
$ws.Range("A62").Value = 45384.440335648149
$ws.Range("B62").Value = "Manya Gupta"
$ws.Range("C62").Value = "B23154"
$ws.Range("D62").Value = "https://www.beecrowd.com.br/judge/en/profile/949111"
$ws.Range("E62").Value = "CSE"
$ws.Range("F62").Value = 0
$ws.Hyperlinks.Add($ws.Range("D62"), "https://www.beecrowd.com.br/judge/en/profile/949111")
$ws.Range("D38").Copy()
$ws.Range("D62").PasteSpecial(-4122)
$ws.Rows.Item(62).RowHeight = 29.4

# Row 63
$ws.Range("A56:F56").Copy()
$ws.Range("A63:F63").PasteSpecial(-4122)
$ws.Range("A63").Value = 45384.446770833332
$ws.Range("B63").Value = "Pranab Ray"
$ws.Range("C63").Value = "B23169"
$ws.Range("D63").Value = "https://www.beecrowd.com.br/judge/en/profile/942674"
$ws.Range("E63").Value = "CSE"
$ws.Range("F63").Value = 0
$ws.Rows.Item(63).RowHeight = 27.6

# Row 64
$ws.Range("A56:F56").Copy()
$ws.Range("A64:F64").PasteSpecial(-4122)
$ws.Range("A64").Value = 45384.453483796293
$ws.Range("B64").Value = "Harshit Kumar Singh"
$ws.Range("C64").Value = "B23133"
$ws.Range("D64").Value = "https://www.beecrowd.com.br/judge/en/profile/949109"
$ws.Range("E64").Value = "CSE"
$ws.Range("F64").Value = 0
$ws.Rows.Item(64).RowHeight = 27.6

# Row 65
$ws.Range("A56:F56").Copy()
$ws.Range("A65:F65").PasteSpecial(-4122)
$ws.Range("A65").Value = 45384.469606481478
$ws.Range("B65").Value = "Shubhankit Singh"
$ws.Range("C65").Value = "B23387"
$ws.Range("D65").Value = "https://www.beecrowd.com.br/judge/en/profile/948383"
$ws.Range("E65").Value = "MSE"
$ws.Range("F65").Value = 0
$ws.Rows.Item(65).RowHeight = 27.6

# Row 66
$ws.Range("A56:F56").Copy()
$ws.Range("A66:F66").PasteSpecial(-4122)
$ws.Range("A66").Value = 45384.486840277779
$ws.Range("B66").Value = "Aaryan Tyagi"
$ws.Range("C66").Value = "B23420"
$ws.Range("D66").Value = "https://www.beecrowd.com.br/judge/en/profile/949129"
$ws.Range("E66").Value = "ME"
$ws.Range("F66").Value = 0
$ws.Hyperlinks.Add($ws.Range("D66"), "https://www.beecrowd.com.br/judge/en/profile/949129")
$ws.Range("D38").Copy()
$ws.Range("D66").PasteSpecial(-4122)
$ws.Rows.Item(66).RowHeight = 29.4

# Row 67
$ws.Range("A56:F56").Copy()
$ws.Range("A67:F67").PasteSpecial(-4122)
$ws.Range("A67").Value = 45384.522835648146
$ws.Range("B67").Value = "Ansh Attre"
$ws.Range("C67").Value = "b23191"
$ws.Range("D67").Value = "https://www.beecrowd.com.br/judge/en/profile/948707"
$ws.Range("E67").Value = "DSE"
$ws.Range("F67").Value = 0
$ws.Rows.Item(67).RowHeight = 27.6

# Row 68
$ws.Range("A56:F56").Copy()
$ws.Range("A68:F68").PasteSpecial(-4122)
$ws.Range("A68").Value = 45384.52925925926
$ws.Range("B68").Value = "Arpita Kumari"
$ws.Range("C68").Value = "B23249"
$ws.Range("D68").Value = "https://www.beecrowd.com.br/judge/en/profile/949137"
$ws.Range("E68").Value = "EE"
$ws.Range("F68").Value = 0
$ws.Rows.Item(68).RowHeight = 27.6

# Row 69
$ws.Range("A56:F56").Copy()
$ws.Range("A69:F69").PasteSpecial(-4122)
$ws.Range("A69").Value = 45384.53224537037
$ws.Range("B69").Value = "Arka"
$ws.Range("C69").Value = "B23120"
$ws.Range("D69").Value = "https://www.beecrowd.com.br/judge/en/profile/948169"
$ws.Range("E69").Value = "CSE"
$ws.Range("F69").Value = 0
$ws.Rows.Item(69).RowHeight = 27.6

# Update the sheet selection to match the new bottom of the data (G69),
# and best-effort scroll so row 56 is near the top of the viewport.
$ws.Activate() | Out-Null
$ws.Range("G69").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 56
$excel.CutCopyMode = $false
